$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking price cells (avoid Excel auto-converting them to Number on write)
$ws.Range("D2:D14").NumberFormat = "@"
$ws.Range("D16:D24").NumberFormat = "@"
$ws.Range("D26:D28").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47:D50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "244.96"

# Row 3
$ws.Range("D3").Value = "25.10"

# Row 4
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "3.499"
$ws.Range("E4").Value = "3LEOLEO"

# Row 5
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "5.018"
$ws.Range("E5").Value = "4HuobiTokenHT"

# Row 6
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "0.05611"
$ws.Range("E6").Value = "5CronosCRO"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "6.574"
$ws.Range("E7").Value = "6KuCoinTokenKCS"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "3.008"
$ws.Range("E8").Value = "7GateTokenGT"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.8113"
$ws.Range("E9").Value = "8MXTokenMX"

# Row 10
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "0.8387"
$ws.Range("E10").Value = "9FTXTokenFTT"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1338"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12
$ws.Range("D12").Value = "0.06948"

# Row 13
$ws.Range("D13").Value = "0.02838"

# Row 14
$ws.Range("D14").Value = "0.09401"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006092"
$ws.Range("E16").Value = "15TigerCashTCH"

# Row 17
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "2.092"
$ws.Range("E17").Value = "16BTSETokenBTSE"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.009672"
$ws.Range("E18").Value = "17OneONEBestin24h"

# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3167"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"

# Row 20
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "0.03276"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"

# Row 21
$ws.Range("D21").Value = "0.1320"

# Row 22
$ws.Range("D22").Value = "3.748"

# Row 23
$ws.Range("D23").Value = "0.04675"

# Row 24
$ws.Range("D24").Value = "0.1369"

# Row 26
$ws.Range("D26").Value = "0.004521"

# Row 27
$ws.Range("D27").Value = "0.00009696"
$ws.Range("E27").Value = "26NitroExNTX"

# Row 28
$ws.Range("D28").Value = "0.0001939"

# Row 40
$ws.Range("D40").Value = "0.03659"

# Row 41
$ws.Range("D41").Value = "0.1349"

# Row 42
$ws.Range("D42").Value = "0.006228"

# Row 43
$ws.Range("D43").Value = "0.002734"

# Row 44
$ws.Range("D44").Value = "0.008086"

# Row 45
$ws.Range("D45").Value = "0.00005289"

# Row 47
$ws.Range("D47").Value = "0.1799"

# Row 48
$ws.Range("D48").Value = "0.002041"

# Row 49
$ws.Range("D49").Value = "0.00002099"

# Row 50
$ws.Range("D50").Value = "0.0001999"
